$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (Initial Weights, Opt Portfolio, Opt Portfolio with View)
$data = @(
    @(2, 0.1,  0.1272849257335129, 0),
    @(3, 0.05, 0.1272908228833229, 0.116373596538509),
    @(4, 0.1,  0.1796008644373881, 0.2176989735329594),
    @(5, 0.1,  0.1722192554824873, 0.1949832312352888),
    @(6, 0.15, 0.138910830905883,  0.07649952086827418),
    @(7, 0.2,  0.1272876666956859, 0.2237036221345661),
    @(8, 0.3,  0.1274056338617202, 0.1707410556904025)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
